$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying shared-string table for "é menor" / "é maior" swaps order
# (see diff). The user-visible effect is that every cell which used to read
# "é menor" now reads "é maior", and every cell which used to read "é maior"
# now reads "é menor". Walk the used range and swap those two values in
# place so the rendered content matches the target workbook.
$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count
$rowOffset = $used.Row
$colOffset = $used.Column

for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
        $text = $cell.Value2
        if ($text -eq "é menor") {
            $cell.Value = "é maior"
        } elseif ($text -eq "é maior") {
            $cell.Value = "é menor"
        }
    }
}
